# UC ELC (transport growth) error fixed
# The sign of the UC_Growth bound (column I) for a set of transport UC
# constraints on the "GrowthConstraintsTransport" sheet was wrong (positive
# instead of negative). Fix the literal entries; the rows below each of them
# hold formulas (=I<row-1>) that simply copy the value down, so they pick the
# correction up automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GrowthConstraintsTransport")

# Rows 30-41 (TFHCV2..5 ELC-N) and 51-62 (TFHCV2..5 HGN-N) each consist of
# 3-row blocks: the first row of every block holds the literal value, the
# following two rows reference the row above with a formula.
$literalRows = 30,33,36,39,51,54,57,60
foreach ($r in $literalRows) {
    $ws.Range("I$r").Value = -2
}

# Move/restore the active selection on this sheet to I1 (matches the
# post-edit saved state).
$ws.Activate()
$ws.Range("I1").Select()
